$d = $word.ActiveDocument

# --- 1. "Lab # 7" -> "Lab # " + "8" (split into two runs, same rPr) ---
# Find the title text and narrow down to just the trailing digit.
$r = $d.Content
$r.Find.Execute("Lab # 7", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$digitStart = $r.End - 1
$digitEnd = $r.End

# Replace the digit's text (briefly leaves a single merged run) ...
$numRange = $d.Range($digitStart, $digitEnd)
$numRange.Text = "8"

# ... then touch formatting on just that character so Word keeps it in
# its own run, matching the committed markup of two adjacent <w:r>
# elements ("Lab # " and "8") that share identical rPr.
$numRange2 = $d.Range($digitStart, $digitStart + 1)
$numRange2.Font.Bold = $true
$numRange2.Font.Bold = $false

# --- 2. Add default/first-page/even-page headers & footers ---
$sec = $d.Sections(1)

# wdHeaderFooterPrimary = 1 -> becomes w:type="default"
$sec.Headers(1).Range.Text = "7/11/2013"

# wdHeaderFooterFirstPage = 2 -> becomes w:type="first"
# wdHeaderFooterEvenPages = 3 -> becomes w:type="even"
# Both are left blank. Writing into their Range is what materializes the
# header1.xml/header2.xml/header3.xml parts (and matching blank footer
# parts for all three variants), without flipping
# DifferentFirstPageHeaderFooter / OddAndEvenPagesHeaderFooter (which
# would additionally stamp a <w:titlePg/> and an evenAndOddHeaders
# setting that this change does not include).
$sec.Headers(2).Range.Text = ""
$sec.Headers(3).Range.Text = ""
